# Fix "Recorded By" (column G) ordering so that dnasr281@gmail.com is
# listed first whenever it is paired with exactly one other recorder.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value()
    if ($null -eq $val) { continue }

    $text = [string]$val
    if ($text -eq "") { continue }

    $parts = $text.Split(",")
    if ($parts.Length -eq 2) {
        $first = $parts[0].Trim()
        $second = $parts[1].Trim()
        if ($second -eq "dnasr281@gmail.com") {
            $cell.Value = "$second, $first"
        }
    }
}
